# Add a small table of numbers to Sheet1 (A1:B4) and leave the
# selection on C1, matching the authored change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = 1
$ws.Range("B1").Value = 8
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 7
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 6
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = 5

$ws.Range("C1").Select() | Out-Null
